# Update the "Please note" wording in the hearing notice paragraph.
# The old wording (spread over three runs: "This case may be released to ",
# "a different court hearing centre", ", in which case you will be notified.")
# is replaced by a longer explanatory sentence that still ends with the
# original wording about the case being released to a different hearing
# centre.

$d = $word.ActiveDocument

$old = "This case may be released to a different court hearing centre, in which case you will be notified."
$new = "Cases are listed in accordance with local hearing arrangements determined by the Judiciary and implemented by the court staff. Every effort is made to ensure that hearings start at the time specified. However, listing practices or other factors may mean that you experience a delay, an adjournment at short notice or your case may be released to a different court hearing centre, in which case you will be notified."

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
